# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
#
# This edit:
#   1. Swaps the data (columns B:AC) between row 89 and row 90 (two match
#      rows whose fixture order was corrected; the "id" in column A stays
#      put, only the underlying match data moves).
#   2. Swaps the data (columns B:AC) between row 117 and row 118, same
#      reasoning as above.
#   3. Appends a brand-new match row (id 142) as row 144, right after the
#      existing last row (143).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap row 89 <-> row 90 (columns B..AC, i.e. 2..29) ---
for ($col = 2; $col -le 29; $col++) {
    $v1 = $ws.Cells.Item(89, $col).Value2
    $v2 = $ws.Cells.Item(90, $col).Value2
    $ws.Cells.Item(89, $col).Value2 = $v2
    $ws.Cells.Item(90, $col).Value2 = $v1
}

# --- 2) Swap row 117 <-> row 118 (columns B..AC, i.e. 2..29) ---
for ($col = 2; $col -le 29; $col++) {
    $v1 = $ws.Cells.Item(117, $col).Value2
    $v2 = $ws.Cells.Item(118, $col).Value2
    $ws.Cells.Item(117, $col).Value2 = $v2
    $ws.Cells.Item(118, $col).Value2 = $v1
}

# --- 3) Append new row 144 ---
$ws.Cells.Item(144, 1).Value2  = "142"
$ws.Cells.Item(144, 2).Value2  = "7862048"
$ws.Cells.Item(144, 3).Value2  = "Lithuania A Lyga"
$ws.Cells.Item(144, 4).Value2  = "Lithuania A Lyga"
$ws.Cells.Item(144, 5).Value2  = "45403.51736111111"
$ws.Cells.Item(144, 6).Value2  = "FK Zalgiris Vilnius"
$ws.Cells.Item(144, 7).Value2  = "FK Dziugas Telsiai"
$ws.Cells.Item(144, 11).Value2 = "1.333"
$ws.Cells.Item(144, 12).Value2 = "5"
$ws.Cells.Item(144, 13).Value2 = "6"
$ws.Cells.Item(144, 14).Value2 = "1.4"
$ws.Cells.Item(144, 15).Value2 = "4.75"
$ws.Cells.Item(144, 16).Value2 = "5.5"
$ws.Cells.Item(144, 17).Value2 = "-1.25"
$ws.Cells.Item(144, 18).Value2 = "1.925"
$ws.Cells.Item(144, 19).Value2 = "1.875"
$ws.Cells.Item(144, 20).Value2 = "2.5"
$ws.Cells.Item(144, 21).Value2 = "1.8"
$ws.Cells.Item(144, 22).Value2 = "2"
$ws.Cells.Item(144, 23).Value2 = "0"
$ws.Cells.Item(144, 24).Value2 = "0"
$ws.Cells.Item(144, 25).Value2 = "0"
$ws.Cells.Item(144, 26).Value2 = "0"
$ws.Cells.Item(144, 27).Value2 = "0"

# Match the formatting (style) of the row above (id column + date column)
$ws.Cells.Item(143, 1).Copy()
$ws.Cells.Item(144, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 5).Copy()
$ws.Cells.Item(144, 5).PasteSpecial(-4122)

Write-Output "Lithuania A Lyga update applied."
